$d = $word.ActiveDocument

function Split-WordAt($findText, $wordStart, $wordLen, $newWord) {
    # Locate the target phrase, then isolate the sub-range that is the
    # misspelled word so it becomes its own run (matching the author's
    # split-run edit) before correcting its text.
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $start = $rng.Start + $wordStart
    $end = $start + $wordLen

    $wordRng = $d.Range($start, $end)
    # Toggling (and restoring) character formatting forces the engine to
    # break this sub-range into its own run, without altering the run's
    # effective formatting.
    $wordRng.Bold = 1
    $wordRng.Bold = 0

    $wordRng = $d.Range($start, $end)
    $wordRng.Text = $newWord
}

# 1) "... Todo este personal esta identificado ..." -> "... está ..."
Split-WordAt "Todo este personal esta identificado" 19 4 "está"

# 2) "Prestigio del médico, se trata de elegir al especialista mas adecuado ..." -> "... más ..."
Split-WordAt "especialista mas adecuado" 13 3 "más"

# 3) "Publicaciones, cuanto mas extensa ..." -> "... más ..."
Split-WordAt "Publicaciones, cuanto mas extensa" 22 3 "más"
